$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Work from the bottom of the document upward so that paragraph indices for
# not-yet-processed regions stay valid while earlier ones shift.
# ---------------------------------------------------------------------------

# --- .post block -----------------------------------------------------------
# Before:
#   26 .post {
#   27   background: #fff;
#   28   margin: 1rem 0;
#   29   padding: 1rem;          (carries a stray lastRenderedPageBreak)
#   30   border-radius: 8px;
#   31   box-shadow: 0 0 5px rgba(0,0,0,0.1);
#   32 }
# After:
#   .post {
#     background: white;
#     padding: 1rem;
#     margin-bottom: 1rem;
#     border-radius: 8px;
#   }

# Remove the box-shadow declaration entirely.
$d.Paragraphs.Item(31).Range.Delete()

# Insert a new "margin-bottom: 1rem;" declaration right after padding.
$d.Paragraphs.Item(29).Range.InsertParagraphAfter()
$d.Paragraphs.Item(30).Range.Text = "  margin-bottom: 1rem;"

# Re-stamp "padding: 1rem;" with identical text so the stale
# lastRenderedPageBreak marker left over from the old layout is dropped.
$d.Paragraphs.Item(29).Range.Text = "  padding: 1rem;"

# Drop the "margin: 1rem 0;" declaration.
$d.Paragraphs.Item(28).Range.Delete()

# background: #fff; -> background: white;
$d.Paragraphs.Item(27).Range.Text = "  background: white;"

# --- nav a block -------------------------------------------------------------
# Before:
#   16 nav a {
#   17   margin: 0 1rem;
#   18   color: #00f0ff;
#   19   text-decoration: none;
#   20 }
# After:
#   nav a {
#     color: #00f0ff;
#     margin: 0 1rem;
#     text-decoration: none;
#   }
#
# The "color: #00f0ff;" paragraph carries spell-check proofErr markup around
# the word "color" - keep that exact paragraph in place and instead relocate
# a freshly-typed "margin: 0 1rem;" paragraph after it, then delete the old
# "margin: 0 1rem;" paragraph. That way the proofErr-bearing paragraph is
# never retyped.
$d.Paragraphs.Item(18).Range.InsertParagraphAfter()
$d.Paragraphs.Item(19).Range.Text = "  margin: 0 1rem;"
$d.Paragraphs.Item(17).Range.Delete()

# --- body block --------------------------------------------------------------
# Before:
#   1 body {
#   2   font-family: Arial, sans-serif;
#   3   background: #f5f7fa;
#   4   color: #333;
#   5   margin: 0;
#   6   padding: 0;
#   7 }
# After:
#   body {
#     font-family: Arial, sans-serif;
#     background: #f5f7fa;
#     margin: 0;
#   }
$d.Paragraphs.Item(6).Range.Delete()
$d.Paragraphs.Item(4).Range.Delete()

Write-Output "done"
